$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.679.35'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '2.349.46'
$ws.Range('E3').Value = '  +3.70%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '233.98'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.79%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.652'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.88%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '66.02'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +4.71%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.456'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.21%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0975'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.52%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '56.42'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '26.94'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +3.07%  '
$ws.Range('D13').Value = '2.696.88'
$ws.Range('E13').Value = '  +3.54%  '
$ws.Range('E14').Value = '  -0.75%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.49'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('E17').Value = '  +1.77%  '
$ws.Range('D18').Value = '2.347.21'
$ws.Range('E18').Value = '  +2.76%  '
$ws.Range('D19').Value = '43.729.66'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').Value = '0.0₃0982'
$ws.Range('E20').Value = '  -1.82%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '74.16'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.90%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.27'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +3.46%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '249.70'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.17%  '
$ws.Range('E24').Value = '  +14.32%  '
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -2.79%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.96'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '22.37'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +7.71%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '174.94'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.92%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.45'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +4.67%  '
$ws.Range('E32').Value = '  -6.11%  '
$ws.Range('E33').Value = '  +1.31%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.02'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +5.14%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0690'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.00'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.29%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.45'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +6.64%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.71'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.41%  '
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0254'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.77%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '9.18'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +11.53%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '18.00'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.04%  '
$ws.Range('E44').Value = '  +10.98%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '99.76'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.98%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0956'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.70%  '
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '4.35'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').Value = '1.448.88'
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.99'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.84%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.32'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.47%  '
